# Daily auto push: insert a new data row for 2026/01/18 at row 647,
# shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 647 (pushes existing row 647 and below down by 1)
$ws.Rows.Item(647).Insert()

# Populate the newly inserted row with the new data point.
# Force column A to be stored as text (not auto-converted to a date serial number).
$ws.Cells.Item(647, 1).NumberFormat = "@"
$ws.Cells.Item(647, 1).Value = "2026/01/18"
$ws.Cells.Item(647, 2).Value = "日"
$ws.Cells.Item(647, 3).Value = 7
$ws.Cells.Item(647, 4).Value = 20
